$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Testcase")

$ws.Range("B2").Value = "Login->To Verify that Successfully landed user role-based landing page, when user Clicking on the Adva pro Login button"
$ws.Range("C2").Value = "openbrowser"
$ws.Range("D2").Value = "https://adva-pro-dev01.paradigmcentral.com"
$ws.Range("E2").Value = "no value"
